$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 68, shifting existing rows 68-70 down to 69-71.
$ws.Rows.Item(68).Insert()

# Copy style (date format) from the row above for column D
$ws.Cells.Item(69, 4).Copy()
$ws.Cells.Item(68, 4).PasteSpecial(-4122)

# Populate the new row 68 with data
$ws.Cells.Item(68, 1).Value = 2
$ws.Cells.Item(68, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(68, 3).Value = "Coquimbo"
$ws.Cells.Item(68, 4).Value = 44595
$ws.Cells.Item(68, 5).Value = 4
$ws.Cells.Item(68, 6).Value = 100112030
$ws.Cells.Item(68, 7).Value = "Poroto granado"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 400
$ws.Cells.Item(68, 11).Value = 20000
$ws.Cells.Item(68, 12).Value = 24000
$ws.Cells.Item(68, 13).Value = 22000
$ws.Cells.Item(68, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(68, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(68, 16).Value = 880
$ws.Cells.Item(68, 17).Value = 25
$ws.Cells.Item(68, 18).Value = "Hortaliza"

$wb.Save()
